$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados Extração VTEX")

$newValue = "American Express, Visa, Mastercard, Boleto Bancário, Hipercard, Elo, Vale, Pix, WH Google Pay, Diners"

foreach ($row in 14..18) {
    $ws.Range("K$row").Value = $newValue
}
